$wb = $excel.ActiveWorkbook

# Rename sheets (task order identifiers refreshed)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16512555831194756"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-1651255585229552"

$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-1651255585230561"

$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-1651255585287304"

$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16512555853492324"

# Sheet 1 (GNG) - update stim filenames
$ws1.Range("B2").Value = "go_stims-16512555830796323.csv"
$ws1.Range("B3").Value = "GNG_stims-16512555831038547.csv"
$ws1.Range("B4").Value = "go_stims-16512555831038547.csv"
$ws1.Range("B5").Value = "GNG_stims-16512555831194756.csv"

# Sheet 2 (NB) - update stim filenames
$ws2.Range("B2").Value = "TB-1651255585028298.csv"
$ws2.Range("B3").Value = "OB-16512555837707882.csv"
$ws2.Range("B4").Value = "ZB-match_0-16512555836697295.csv"
$ws2.Range("B5").Value = "OB-16512555839997113.csv"
$ws2.Range("B6").Value = "TB-16512555852115536.csv"
$ws2.Range("B7").Value = "ZB-match_2-1651255583458377.csv"
$ws2.Range("B8").Value = "TB-16512555846689126.csv"
$ws2.Range("B9").Value = "ZB-match_6-16512555833310513.csv"
$ws2.Range("B10").Value = "OB-1651255584419442.csv"

# Sheet 3 (RS) - no content changes

# Sheet 4 (TOL) - update stim filenames
$ws4.Range("B2").Value = "MM_stims-16512555852493315.csv"
$ws4.Range("B3").Value = "ZM_stims-16512555852325563.csv"
$ws4.Range("B4").Value = "MM_stims-16512555852708528.csv"
$ws4.Range("B5").Value = "ZM_stims-16512555852493315.csv"
$ws4.Range("B6").Value = "MM_stims-16512555852863033.csv"
$ws4.Range("B7").Value = "ZM_stims-16512555852708528.csv"

# Sheet 5 (vSAT) - update stim filenames
$ws5.Range("B2").Value = "SAT_stims-16512555852903044.csv"
$ws5.Range("B3").Value = "vSAT_stims-1651255585319236.csv"
$ws5.Range("B4").Value = "vSAT_stims-1651255585333234.csv"
$ws5.Range("B5").Value = "SAT_stims-16512555853022342.csv"
